# Update the "Sleep Disorders" effects table with refreshed model output.
# Because several cells in the IRR column share identical text (e.g. the
# Canada/China/Norway "Level change" and "Trend change" rows were all
# "1.19 (1.16 to 1.22), p < 0.001" / "1.01 (1 to 1.01), p < 0.001"),
# Find/Replace would affect the wrong occurrence. Instead we address each
# target cell directly via the Tables/Cell object model and overwrite its
# Range.Text, which only ever touches that single cell's run.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Table layout (row 1 is the header row):
#   Row 10 -> Norway / Level change
#   Row 11 -> Norway / Trend change
#   Row 16 -> USA    / Level change
#   Row 17 -> USA    / Trend change

$t.Cell(10, 3).Range.Text = "1 (0.94 to 1.06), p = 0.989"
$t.Cell(11, 3).Range.Text = "1 (0.99 to 1.01), p = 0.894"
$t.Cell(16, 3).Range.Text = "1.19 (1.03 to 1.36), p = 0.027"
$t.Cell(17, 3).Range.Text = "1.02 (1.01 to 1.03), p < 0.001"
